$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 124, shifting existing rows 124-149 down to 125-150.
$ws.Rows.Item(124).Insert()

# Populate the new row 124 with the "Cacho cabra verde" record.
$ws.Range("A124").Value = 7
$ws.Range("B124").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C124").Value = "Ñuble"
$ws.Range("D124").Value = 44995
$ws.Range("E124").Value = 16
$ws.Range("F124").Value = 100112021
$ws.Range("G124").Value = "Ají"
$ws.Range("H124").Value = "Cacho cabra verde"
$ws.Range("I124").Value = "Primera"
$ws.Range("J124").Value = 30
$ws.Range("K124").Value = 14000
$ws.Range("L124").Value = 14000
$ws.Range("M124").Value = 14000
$ws.Range("N124").Value = "$/saco 25 kilos"
$ws.Range("O124").Value = "Región del Maule"
$ws.Range("P124").Value = 560
$ws.Range("Q124").Value = 25
$ws.Range("R124").Value = "Hortaliza"
